# productVerify and BORetail updated: stale element exception handled
#
# Appends "1" to the test usernames (column A) and to the local-part of the
# corresponding test email addresses (column C) for data rows 2-6, then
# moves the active selection to A6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ahostess): update username then email together.
$a2 = $ws.Cells.Item(2, 1)
$c2 = $ws.Cells.Item(2, 3)
$a2.Value = $a2.Value() + "1"
$c2.Value = ($c2.Value() -replace '@test\.com$', '1@test.com')

# Update the remaining emails (rows 3-6) first...
for ($row = 3; $row -le 6; $row++) {
    $emailCell = $ws.Cells.Item($row, 3)  # column C
    $emailCell.Value = ($emailCell.Value() -replace '@test\.com$', '1@test.com')
}

# ...then update the remaining usernames (rows 3-6).
for ($row = 3; $row -le 6; $row++) {
    $nameCell = $ws.Cells.Item($row, 1)  # column A
    $nameCell.Value = $nameCell.Value() + "1"
}

$ws.Range("A6").Select()
